$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 6.206015333333333
$ws.Range("H2").Value2 = 18.618046
$ws.Range("I2").Value2 = 0.0150172404156507
$ws.Range("J2").Value2 = 0.0150172404156507
$ws.Range("M2").Value2 = 0.8366046666666667
$ws.Range("N2").Value2 = 2.509814
$ws.Range("O2").Value2 = 0.08025679986157715
$ws.Range("P2").Value2 = 0.08025679986157715
$ws.Range("Q2").Value2 = 5.191981389271556
$ws.Range("R2").Value2 = 46.727832503444
$ws.Range("S2").Value2 = 0.001205235658512066
$ws.Range("T2").Value2 = 0.001205235658512066
$ws.Range("G3").Value2 = 6.206015333333333
$ws.Range("H3").Value2 = 18.618046
$ws.Range("I3").Value2 = 0.0150172404156507
$ws.Range("J3").Value2 = 0.0150172404156507
$ws.Range("M3").Value2 = 7.939250333333333
$ws.Range("O3").Value2 = 0.7616247559221037
$ws.Range("P3").Value2 = 0.7616247559221038
$ws.Range("Q3").Value2 = 49.27110930383844
$ws.Range("R3").Value2 = 443.4399837345459
$ws.Range("S3").Value2 = 0.01143750206619351
$ws.Range("T3").Value2 = 0.01143750206619351
$ws.Range("G4").Value2 = 6.206015333333333
$ws.Range("H4").Value2 = 18.618046
$ws.Range("I4").Value2 = 0.0150172404156507
$ws.Range("J4").Value2 = 0.0150172404156507
$ws.Range("M4").Value2 = 1.648242
$ws.Range("N4").Value2 = 4.944726
$ws.Range("O4").Value2 = 0.1581184442163192
$ws.Range("P4").Value2 = 0.1581184442163192
$ws.Range("Q4").Value2 = 10.229015125044
$ws.Range("R4").Value2 = 92.061136125396
$ws.Range("S4").Value2 = 0.002374502690945119
$ws.Range("T4").Value2 = 0.002374502690945119
$ws.Range("I5").Value2 = 0.9317452840597572
$ws.Range("J5").Value2 = 0.9317452840597571
$ws.Range("M5").Value2 = 0.8366046666666667
$ws.Range("N5").Value2 = 2.509814
$ws.Range("O5").Value2 = 0.08025679986157715
$ws.Range("P5").Value2 = 0.08025679986157715
$ws.Range("Q5").Value2 = 322.1366935923949
$ws.Range("R5").Value2 = 2899.230242331554
$ws.Range("S5").Value2 = 0.07477889478475228
$ws.Range("T5").Value2 = 0.07477889478475228
$ws.Range("I6").Value2 = 0.9317452840597572
$ws.Range("J6").Value2 = 0.9317452840597571
$ws.Range("M6").Value2 = 7.939250333333333
$ws.Range("O6").Value2 = 0.7616247559221037
$ws.Range("P6").Value2 = 0.7616247559221038
$ws.Range("R6").Value2 = 27513.25158100266
$ws.Range("S6").Value2 = 0.7096402745535837
$ws.Range("T6").Value2 = 0.7096402745535837
$ws.Range("I7").Value2 = 0.9317452840597572
$ws.Range("J7").Value2 = 0.9317452840597571
$ws.Range("M7").Value2 = 1.648242
$ws.Range("N7").Value2 = 4.944726
$ws.Range("O7").Value2 = 0.1581184442163192
$ws.Range("P7").Value2 = 0.1581184442163192
$ws.Range("Q7").Value2 = 634.6596538071541
$ws.Range("R7").Value2 = 5711.936884264387
$ws.Range("S7").Value2 = 0.1473261147214212
$ws.Range("T7").Value2 = 0.1473261147214212
$ws.Range("G8").Value2 = 22.00088566666667
$ws.Range("H8").Value2 = 66.002657
$ws.Range("I8").Value2 = 0.05323747552459213
$ws.Range("J8").Value2 = 0.05323747552459213
$ws.Range("M8").Value2 = 0.8366046666666667
$ws.Range("N8").Value2 = 2.509814
$ws.Range("O8").Value2 = 0.08025679986157715
$ws.Range("P8").Value2 = 0.08025679986157715
$ws.Range("Q8").Value2 = 18.40604361953311
$ws.Range("R8").Value2 = 165.654392575798
$ws.Range("S8").Value2 = 0.004272669418312803
$ws.Range("T8").Value2 = 0.004272669418312802
$ws.Range("G9").Value2 = 22.00088566666667
$ws.Range("H9").Value2 = 66.002657
$ws.Range("I9").Value2 = 0.05323747552459213
$ws.Range("J9").Value2 = 0.05323747552459213
$ws.Range("M9").Value2 = 7.939250333333333
$ws.Range("O9").Value2 = 0.7616247559221037
$ws.Range("P9").Value2 = 0.7616247559221038
$ws.Range("Q9").Value2 = 174.6705388627119
$ws.Range("R9").Value2 = 1572.034849764407
$ws.Range("S9").Value2 = 0.04054697930232645
$ws.Range("T9").Value2 = 0.04054697930232645
$ws.Range("G10").Value2 = 22.00088566666667
$ws.Range("H10").Value2 = 66.002657
$ws.Range("I10").Value2 = 0.05323747552459213
$ws.Range("J10").Value2 = 0.05323747552459213
$ws.Range("M10").Value2 = 1.648242
$ws.Range("N10").Value2 = 4.944726
$ws.Range("O10").Value2 = 0.1581184442163192
$ws.Range("P10").Value2 = 0.1581184442163192
$ws.Range("Q10").Value2 = 36.262783792998
$ws.Range("R10").Value2 = 326.365054136982
$ws.Range("S10").Value2 = 0.00841782680395288
$ws.Range("T10").Value2 = 0.00841782680395288
